$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Setup")

# --- Row 20: AXDWAdmin -----------------------------------------------------
$ws.Range("A20").Value = "AXDWAdmin"
$ws.Range("B20").Formula = "=LEFT(A20,4) & `$G`$2 & MID(A20,5,10)"
$ws.Range("C20").Formula = "= LEFT(A20,4) & `$G`$2 & MID(A20,5,10)"
$ws.Range("D20").Value = $ws.Range("D18").Value2

# --- Row 21: AXDWRuntimeuser ------------------------------------------------
$ws.Range("A21").Value = "AXDWRuntimeuser"
$ws.Range("B21").Formula = "=LEFT(A21,4) & `$G`$2 & MID(A21,5,11)"
$ws.Range("C21").Formula = "= LEFT(A21,4) & `$G`$2 & MID(A21,5,11)"
$ws.Range("D21").Value = $ws.Range("D18").Value2

# New cellXfs style (numFmtId 0 / General, with applyNumberFormat flag set)
# gets minted the same way Excel mints it when the Number Format is
# explicitly (re)applied to the "Account Name" cells from the UI.
$ws.Range("B20:B21").NumberFormat = "General"

# --- Extend the "Table33" structured table to include the two new rows -----
$lo = $ws.ListObjects.Item("Table33")
$lo.Resize($ws.Range("A12:D21"))

# --- Selection, matching where the author ended up after the edit ---------
$ws.Range("H20").Select()
